$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the existing headers (bold, bordered, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill column I with constant 1 for all data rows
$ws.Range("I2:I38").Value = 1

# Column J mirrors column H for all data rows
$ws.Range("H2:H38").Copy()
$ws.Range("J2:J38").PasteSpecial(-4163)

$excel.CutCopyMode = 0
